# Update "want to go" counts (column F) on both the "展览" sheet and the
# "全部类型" sheet, which mirror the same rows.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8373
    3  = 7908
    10 = 176
    12 = 715
    13 = 133
    16 = 56
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
